# Generate Report for Handback
# - Marks a.md as "handed back" (status text) across all sheets
# - Fills in the "Latest Target File" / "Latest Handback File" columns
#   (F/G) for the zh-cn and de-de report sheets, with hyperlinks
# - Stamps the "Latest Handback DateTime" (H) for each language

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status column updates -------------------------------------------------
# These all point at the same shared string, so update every occurrence to
# the new text so the whole workbook stays in sync.
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# --- zh-cn: Latest Target File (F) / Latest Handback File (G) --------------
$zhMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/15fce05df4cae3230b019c592856b01ddb31362b/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2195455a5e7089027f6bc0bdbee19e93532cdc4e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhMdUrl, "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhMdUrl, "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

# zh-cn: Latest Handback DateTime (H) - both rows converge on the same
# stamp, so they stay on one shared string after save.
$wsZhCn.Range("H2").Value = "2016-03-23 10:35:38"
$wsZhCn.Range("H3").Value = "2016-03-23 10:35:38"

# --- de-de: Latest Target File (F) / Latest Handback File (G) --------------
$deMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/15fce05df4cae3230b019c592856b01ddb31362b/e2e/a.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea321dd325f0d742f0d21f8fe6f838a6f6d99a98/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deMdUrl, "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deMdUrl, "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deXlfUrl, "", "", $deXlfName)

# de-de: Latest Handback DateTime (H) - a different stamp than zh-cn, so it
# lands on its own shared string.
$wsDeDe.Range("H2").Value = "2016-03-23 10:35:47"
$wsDeDe.Range("H3").Value = "2016-03-23 10:35:47"

Write-Output "Handback report generated"
